# Adds all IG authors as contact
# - Duplicates the existing "Contact" row (row 11) twice, inserting two new
#   rows right after it, so the Metadata sheet gains two more Contact rows.
# - Bumps the "Date" metadata value to the new export timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new blank rows right after the existing "Contact" row (row 11),
# pushing Jurisdiction/Description/... etc. down by two rows.
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(12).Insert()

# Copy the formatting of the original Contact row (11) onto the two new rows
# (12 and 13) so they keep the same style as the rest of the data rows.
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B13").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new Contact rows with the same content as the existing one.
$ws.Range("A12").Value = "Contact"
$ws.Range("B12").Value = "No display for ContactDetail"
$ws.Range("A13").Value = "Contact"
$ws.Range("B13").Value = "No display for ContactDetail"

# Update the Date metadata value to reflect the new export time.
$ws.Range("B8").Value = "2022-01-21T07:49:24+01:00"
